$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SheetName1")

# Update T2 value
$ws.Range("T2").Value = 445300

# Update selection to T2
$ws.Range("T2").Select()
